$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timesheet row bookkeeping -----------------------------------------
# C7 ("Total Horas Diarias" for 21/07) was never filled in - add the
# missing 1 hour total (same time format as the rest of column C).
$ws.Range("C7").Value = 1/24
$ws.Range("C7").NumberFormat = "h:mm"

# The day row holding 22/07/2014 (serial 41842) is removed entirely,
# shifting the remaining date rows up by one (23/07 -> row 8, 24/07 ->
# row 9, 25/07 -> row 10).
$ws.Rows(8).Delete()

# --- Novos calculos auxiliares ------------------------------------------
# Calculo do volume da encomenda: 25 caixas * 400 unidades.
$ws.Range("F17").Formula = "=25*400"

# Match the selection left behind by the author after these edits.
$ws.Range("C18").Select()
